# Applies the "Updated symbol list" data refresh to cryptos.xlsx (Sheet1).
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). F/G (date/hour) are left as-is.
# NumberFormat is forced to "@" (Text) before each write so that numeric-looking
# strings (e.g. "247.03", "0.83%") are stored as literal text, matching the
# original inlineStr text cells instead of being auto-coerced into Number cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '247.03'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.83%'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.42'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '5.22%'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.082'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.74%'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05602'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.19%'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.476'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-1.49%'

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8133'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.55%'

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8445'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '0.25%'

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07004'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.86%'

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.02853'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.48%'

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09399'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.03%'

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.001509'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.97%'

# Row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'One'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0005990'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.37%'

# Row 14
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'TigerCash'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.006181'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.33%'

# Row 15
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'LEO'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.607'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '3.15%'

# Row 16
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'GateToken'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.015'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.23%'

# Row 17
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'BTSEToken'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.056'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.72%'

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.3126'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-2.21%'

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.1342'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.36%'

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03178'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-2.69%'

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.45%'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.743'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '0.14%'

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04648'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.59%'

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-1.46%'

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001244'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.10%'

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004591'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '1.44%'

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.00009600'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-0.99%'

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001398'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '-27.94%'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03669'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.51%'

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1351'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.75%'

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.002660'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-2.73%'

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003421'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-45.07%'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008945'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '10.92%'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005360'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '1.40%'

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.02%'

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-38.87%'

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002609'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '27.80%'

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.02%'

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.02%'
